$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update formula/label text in columns B, C, D for rows 2-5 ---
$ws.Range("B2").Value = "TH - (0.001493 * IC)"
$ws.Range("C2").Value = "495 /eBD"
$ws.Range("D2").Value = "metoda0"

$ws.Range("C3").Value = "0.21 * (AB + TR+ SS)"
$ws.Range("D3").Value = "metoda1"

$ws.Range("B4").Value = "1.1145 - (0.000465 * AY)"
$ws.Range("C4").Value = "4.570/ eBD - 4.142"
$ws.Range("D4").Value = "metoda2"

$ws.Range("C5").Value = "(15 * BM) + 8 * (AB + TR+ SS)"
$ws.Range("D5").Value = "metoda3"

# --- Row 4 / row 5 no longer need the taller custom row height ---
$ws.Rows(4).AutoFit()
$ws.Rows(5).AutoFit()

# --- B5 becomes an (empty) red Arial 10 wrap-text cell ---
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 10
$ws.Range("B5").Font.Color = 255
$ws.Range("B5").WrapText = $true

# --- B3 becomes an (empty) plain cell ---
$ws.Range("B3").WrapText = $false

# --- Selection moves to C5 ---
$ws.Range("C5").Select()
